$wb = $excel.ActiveWorkbook

# --- Rename the "Include ValueSets*" sheets to "Include ValueSet #N" ---
$wb.Worksheets.Item(2).Name = "Include ValueSet #0"
$wb.Worksheets.Item(3).Name = "Include ValueSet #1"
$wb.Worksheets.Item(4).Name = "Include ValueSet #2"
$wb.Worksheets.Item(5).Name = "Include ValueSet #3"
$wb.Worksheets.Item(6).Name = "Include ValueSet #4"

# --- Update the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item(1)

# Version: 0.2.1 -> 0.2.2
$wsMeta.Range("B3").Value = "0.2.2"

# Date: refreshed timestamp
$wsMeta.Range("B8").Value = "2024-09-11T16:17:59-05:00"

# Contact: now shows the publisher contact string
$wsMeta.Range("B10").Value = "MITRE (https://www.mitre.org)"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# before "Description" (old row 11), copying formatting from the row below
# so the new row keeps the same cell style as the rest of the table.
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
